$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows shift by 0 relative to data rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 797
$ws1.Range("F6").Value  = 94
$ws1.Range("F7").Value  = 322
$ws1.Range("F8").Value  = 4050
$ws1.Range("F10").Value = 4760
$ws1.Range("F11").Value = 527
$ws1.Range("F12").Value = 1199
$ws1.Range("F13").Value = 82

# Sheet "全部类型" (All types) - has one extra row vs "展览" so the same
# events land one row further down starting at row 8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 797
$ws4.Range("F6").Value  = 94
$ws4.Range("F8").Value  = 322
$ws4.Range("F9").Value  = 4050
$ws4.Range("F11").Value = 4760
$ws4.Range("F12").Value = 527
$ws4.Range("F13").Value = 1199
$ws4.Range("F14").Value = 82
